$d = $word.ActiveDocument

# 1. Expand the "Last contract finished..." sentence with the Morocco story.
$d.Content.Find.Execute(
    "Last contract finished on the 1st of January 2020.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Last contract finished on the 1st of January 2020. I then went to Marrakech, Morocco on the 4th of January 2020 to buy an apartment, which took around 3 months. By the time I had completed the whole process it was the middle of March 2020 and the borders in Morocco were closed so I could not return to London. I returned to London on the 22nd of July 2020, after the Moroccan government opened the borders on the 15th of July 2020.",
    2) | Out-Null

# 2. Locate the "Kind regards," / "Mohamed Bana" paragraphs that currently sit
#    right after the paragraph above (they need to move to the very end of the
#    letter, after the horizontal-rule separator, replacing the
#    website/e-mail/phone line there).
$signOffPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Kind regards,") {
        $signOffPara = $i
        break
    }
}

$kindRegardsRange = $d.Paragraphs.Item($signOffPara).Range
$mohamedBanaRange = $d.Paragraphs.Item($signOffPara + 1).Range

# Delete "Mohamed Bana" paragraph (incl. its own mark) then "Kind regards,"
# paragraph (incl. its own mark) - removing both from their old location.
$mohamedBanaRange.Delete() | Out-Null
$kindRegardsRange.Delete() | Out-Null

# 3. Find the contact-details paragraph (the one with the bana.io hyperlinks)
#    that follows the horizontal rule, and swap its contents for "Kind
#    regards," / "Mohamed Bana" while keeping its own (FirstParagraph) style.
$contactPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*bana.io*") {
        $contactPara = $i
        break
    }
}

$contactRange = $d.Paragraphs.Item($contactPara).Range
$contentOnly = $d.Range($contactRange.Start, $contactRange.End - 1)
$contentOnly.Delete() | Out-Null

$contactPara = $d.Paragraphs.Item($contactPara)
$insertPoint = $d.Range($contactPara.Range.Start, $contactPara.Range.Start)
$insertPoint.InsertAfter("Kind regards,") | Out-Null

$kindRegardsPara = $d.Paragraphs.Item($contactPara.Index)
$kindRegardsPara.Range.InsertParagraphAfter() | Out-Null

$mohamedPara = $d.Paragraphs.Item($kindRegardsPara.Index + 1)
$mohamedPara.Style = "Body Text"
$mohamedInsertPoint = $d.Range($mohamedPara.Range.Start, $mohamedPara.Range.Start)
$mohamedInsertPoint.InsertAfter("Mohamed Bana") | Out-Null
